# Update countries & provincias Spain
# - Reorder "Somalia" ahead of "Lituania"/"Mayotte" in the country list and
#   refresh its case numbers. Lituania and Mayotte shift down one row each,
#   keeping their own previous numbers, and the row that used to hold
#   Somalia (now occupied by Mayotte) disappears as a distinct "Somalia"
#   entry.
# - Refresh COVID-19 counters for several countries (Estados Unidos,
#   Singapur, Republica Dominicana, Irak, Grecia, Somalia, Lituania,
#   Mayotte, Birmania).
# - Bump the "Datos actualizados" timestamp string from 17:05 to 17:35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp footer (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 17:35"

# --- Helper: write a full data row (Pais, Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ----
function Set-Row($Row, $Pais, $CasosTotales, $NuevosCasos, $CasosActivos, $Recuperados, $CasosCriticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Pais
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $CasosCriticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Estados Unidos (row 4)
Set-Row 4 "Estados Unidos" 1691225 4789 451776 1140053 0 96 99396

# Singapur (row 29)
Set-Row 29 "Singapur" 31960 344 15738 16199 0 0 23

# Republica Dominicana (row 45)
Set-Row 45 "Republica Dominicana" 15073 272 8285 6328 0 2 460

# Irak (row 69)
Set-Row 69 "Irak" 4632 163 2811 1658 0 3 163

# Grecia (row 81)
Set-Row 81 "Grecia" 2882 4 1374 1336 0 1 172

# Somalia moves up (was listed after Mayotte, now right after Islandia) and
# gets fresh numbers; Lituania and Mayotte shift down one row each, keeping
# their previous numbers.
Set-Row 94 "Somalia"  1689 95 235  1388 0 5 66
Set-Row 95 "Lituania" 1635 12 1138 434  0 0 63
Set-Row 96 "Mayotte"  1609 22 894  695  0 0 20

# Birmania (row 156)
Set-Row 156 "Birmania" 203 2 123 74 0 0 6
